# Bond dates update: one day has passed, so for every bond row:
#   - column G ("Dni od poprzedniej wypłaty" / days since previous payout) increases by 1
#   - column I ("Dni do następnej wypłaty" / days until next payout) decreases by 1
# Only rows/cells that already contain a numeric value are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$colG = 7   # column G
$colI = 9   # column I

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, $colG)
    $gVal = $gCell.Value()
    if ($gVal -ne $null) {
        $gCell.Value = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, $colI)
    $iVal = $iCell.Value()
    if ($iVal -ne $null) {
        $iCell.Value = $iVal - 1
    }
}
